# Fix missing "Retiring by 2050" capacity data (column C) that was
# incorrectly left at 0 for several countries/regions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11  = 3.236
    19  = 6.495
    43  = 0.345
    59  = 0.11
    103 = 3.068
    131 = 30.0296
    135 = 3.9876
    151 = 33.3296
    187 = 0.4
    247 = 3
    315 = 23.1346
    343 = 3.068
    375 = 6.15
    379 = 12.843
    387 = 0.3
    435 = 3.8776
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
